$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Cells.Item(2, 2).Value = 19.3032722193255
$ws.Cells.Item(2, 4).Value = 8.618142845807071
$ws.Cells.Item(2, 5).Value = 14.73139184132897
$ws.Cells.Item(2, 6).Value = 39.46283821389006
$ws.Cells.Item(2, 7).Value = 47.82749118213547
$ws.Cells.Item(2, 8).Value = 18.55754805184232
$ws.Cells.Item(2, 10).Value = 11.0317279510894
$ws.Cells.Item(2, 12).Value = 11.47213404137494
$ws.Cells.Item(2, 13).Value = 18.03009658347736
$ws.Cells.Item(2, 14).Value = 20.21711124810983

$ws.Cells.Item(3, 2).Value = 19.10327141264699
$ws.Cells.Item(3, 4).Value = 8.533844714924252
$ws.Cells.Item(3, 5).Value = 14.57916417425866
$ws.Cells.Item(3, 6).Value = 39.42211884920286
$ws.Cells.Item(3, 7).Value = 47.43139997684967
$ws.Cells.Item(3, 8).Value = 18.55971654144591
$ws.Cells.Item(3, 10).Value = 10.97962870786365
$ws.Cells.Item(3, 12).Value = 11.33607853651253
$ws.Cells.Item(3, 13).Value = 17.91394446464671
$ws.Cells.Item(3, 14).Value = 20.28241678658563

$ws.Cells.Item(4, 2).Value = 18.98322673041606
$ws.Cells.Item(4, 4).Value = 8.480869230386794
$ws.Cells.Item(4, 5).Value = 14.48384913065668
$ws.Cells.Item(4, 6).Value = 39.40878372390222
$ws.Cells.Item(4, 7).Value = 47.20425514132499
$ws.Cells.Item(4, 8).Value = 18.56573348416452
$ws.Cells.Item(4, 10).Value = 10.94724088238149
$ws.Cells.Item(4, 12).Value = 11.25362415191868
$ws.Cells.Item(4, 13).Value = 17.84526161643937
$ws.Cells.Item(4, 14).Value = 20.3246940080631

$ws.Cells.Item(5, 2).Value = 18.93505064918466
$ws.Cells.Item(5, 4).Value = 8.458982493054945
$ws.Cells.Item(5, 5).Value = 14.44456053826033
$ws.Cells.Item(5, 6).Value = 39.40628142813265
$ws.Cells.Item(5, 7).Value = 47.11581848881388
$ws.Cells.Item(5, 8).Value = 18.56936157324085
$ws.Cells.Item(5, 10).Value = 10.93394560308952
$ws.Cells.Item(5, 12).Value = 11.22033040248681
$ws.Cells.Item(5, 13).Value = 17.81795748748227
$ws.Cells.Item(5, 14).Value = 20.34247182625404

$ws.Cells.Item(6, 2).Value = 18.92709736019268
$ws.Cells.Item(6, 4).Value = 8.455330376537804
$ws.Cells.Item(6, 5).Value = 14.43801021634751
$ws.Cells.Item(6, 6).Value = 39.40604289663667
$ws.Cells.Item(6, 7).Value = 47.10138531093166
$ws.Cells.Item(6, 8).Value = 18.57003497652105
$ws.Cells.Item(6, 10).Value = 10.9317321699246
$ws.Cells.Item(6, 12).Value = 11.21482154359892
$ws.Cells.Item(6, 13).Value = 17.81346563622271
$ws.Cells.Item(6, 14).Value = 20.34545705241135

$ws.Cells.Item(7, 2).Value = 18.9825739373057
$ws.Cells.Item(7, 4).Value = 8.480575258659623
$ws.Cells.Item(7, 5).Value = 14.48332105738256
$ws.Cells.Item(7, 6).Value = 39.40873811056228
$ws.Cells.Item(7, 7).Value = 47.20304563428395
$ws.Cells.Item(7, 8).Value = 18.56577765505148
$ws.Cells.Item(7, 10).Value = 10.94706196593159
$ws.Cells.Item(7, 12).Value = 11.25317385233755
$ws.Cells.Item(7, 13).Value = 17.84489058277393
$ws.Cells.Item(7, 14).Value = 20.32493153904376

$ws.Cells.Item(8, 2).Value = 19.23377016527638
$ws.Cells.Item(8, 4).Value = 8.589333412194931
$ws.Cells.Item(8, 5).Value = 14.67929737883385
$ws.Cells.Item(8, 6).Value = 39.44637644049976
$ws.Cells.Item(8, 7).Value = 47.68764267172041
$ws.Cells.Item(8, 8).Value = 18.5573222323046
$ws.Cells.Item(8, 10).Value = 11.01384884386907
$ws.Cells.Item(8, 12).Value = 11.42501787450292
$ws.Cells.Item(8, 13).Value = 17.98951399373777
$ws.Cells.Item(8, 14).Value = 20.23917716687568

$ws.Cells.Item(9, 2).Value = 19.74598579826261
$ws.Cells.Item(9, 4).Value = 8.792659957724654
$ws.Cells.Item(9, 5).Value = 15.04820675313479
$ws.Cells.Item(9, 6).Value = 39.61272044455223
$ws.Cells.Item(9, 7).Value = 48.76097908927409
$ws.Cells.Item(9, 8).Value = 18.57798254158366
$ws.Cells.Item(9, 10).Value = 11.14151242882574
$ws.Cells.Item(9, 12).Value = 11.76907516915124
$ws.Cells.Item(9, 13).Value = 18.29301250531031
$ws.Cells.Item(9, 14).Value = 20.08823788825997

$ws.Cells.Item(10, 2).Value = 20.13119727828869
$ws.Cells.Item(10, 4).Value = 8.935592981367463
$ws.Cells.Item(10, 5).Value = 15.30888368153649
$ws.Cells.Item(10, 6).Value = 39.79113026620357
$ws.Cells.Item(10, 7).Value = 49.61816240011567
$ws.Cells.Item(10, 8).Value = 18.61591136981767
$ws.Cells.Item(10, 10).Value = 11.23310465869188
$ws.Cells.Item(10, 12).Value = 12.02410509667364
$ws.Cells.Item(10, 13).Value = 18.52670289815836
$ws.Cells.Item(10, 14).Value = 19.9877499161707

$ws.Cells.Item(11, 2).Value = 20.30770805187069
$ws.Cells.Item(11, 4).Value = 8.999137064259974
$ws.Cells.Item(11, 5).Value = 15.42502341362287
$ws.Cells.Item(11, 6).Value = 39.88439273382043
$ws.Cells.Item(11, 7).Value = 50.02151009768458
$ws.Cells.Item(11, 8).Value = 18.63810075306854
$ws.Cells.Item(11, 10).Value = 11.274254125111
$ws.Cells.Item(11, 12).Value = 12.14017909429509
$ws.Cells.Item(11, 13).Value = 18.63502389100847
$ws.Cells.Item(11, 14).Value = 19.94427557853449

$ws.Cells.Item(12, 2).Value = 20.37467626812819
$ws.Cells.Item(12, 4).Value = 9.022980610210162
$ws.Cells.Item(12, 5).Value = 15.46863476556154
$ws.Cells.Item(12, 6).Value = 39.92143584977767
$ws.Cells.Item(12, 7).Value = 50.1760422665071
$ws.Cells.Item(12, 8).Value = 18.64721118568637
$ws.Cells.Item(12, 10).Value = 11.28975860246421
$ws.Cells.Item(12, 12).Value = 12.184106098241
$ws.Cells.Item(12, 13).Value = 18.67630437624582
$ws.Cells.Item(12, 14).Value = 19.9281333796101

$ws.Cells.Item(13, 2).Value = 20.36024868447692
$ws.Cells.Item(13, 4).Value = 9.017855332863467
$ws.Cells.Item(13, 5).Value = 15.45925892190086
$ws.Cells.Item(13, 6).Value = 39.91338140698075
$ws.Cells.Item(13, 7).Value = 50.14268343996511
$ws.Cells.Item(13, 8).Value = 18.64521765021441
$ws.Cells.Item(13, 10).Value = 11.28642296069058
$ws.Cells.Item(13, 12).Value = 12.1746474334287
$ws.Cells.Item(13, 13).Value = 18.66740269566558
$ws.Cells.Item(13, 14).Value = 19.93159565233096

$ws.Cells.Item(14, 2).Value = 20.3132153329307
$ws.Cells.Item(14, 4).Value = 9.001103119170683
$ws.Cells.Item(14, 5).Value = 15.42861882087873
$ws.Cells.Item(14, 6).Value = 39.88740577723407
$ws.Cells.Item(14, 7).Value = 50.03418840247866
$ws.Cells.Item(14, 8).Value = 18.63883610225646
$ws.Cells.Item(14, 10).Value = 11.27553127081394
$ws.Cells.Item(14, 12).Value = 12.14379375558517
$ws.Cells.Item(14, 13).Value = 18.63841497831025
$ws.Cells.Item(14, 14).Value = 19.94294113211025

$ws.Cells.Item(15, 2).Value = 20.28442098544346
$ws.Cells.Item(15, 4).Value = 8.990813155717277
$ws.Cells.Item(15, 5).Value = 15.40980242550845
$ws.Cells.Item(15, 6).Value = 39.87171935916216
$ws.Cells.Item(15, 7).Value = 49.96796150312727
$ws.Cells.Item(15, 8).Value = 18.63501932584185
$ws.Cells.Item(15, 10).Value = 11.26884953763581
$ws.Cells.Item(15, 12).Value = 12.12489031588082
$ws.Cells.Item(15, 13).Value = 18.62069240691746
$ws.Cells.Item(15, 14).Value = 19.94993227952259

$ws.Cells.Item(16, 2).Value = 20.11968275665822
$ws.Cells.Item(16, 4).Value = 8.931409918248649
$ws.Cells.Item(16, 5).Value = 15.30124305436155
$ws.Cells.Item(16, 6).Value = 39.78527777467382
$ws.Cells.Item(16, 7).Value = 49.59206105825405
$ws.Cells.Item(16, 8).Value = 18.61456040852352
$ws.Cells.Item(16, 10).Value = 11.23040473972561
$ws.Cells.Item(16, 12).Value = 12.01651747042993
$ws.Cells.Item(16, 13).Value = 18.5196620145437
$ws.Cells.Item(16, 14).Value = 19.99063599580353

$ws.Cells.Item(17, 2).Value = 20.01890797659304
$ws.Cells.Item(17, 4).Value = 8.894584765581989
$ws.Cells.Item(17, 5).Value = 15.2340074034585
$ws.Cells.Item(17, 6).Value = 39.73533961323741
$ws.Cells.Item(17, 7).Value = 49.36479809548497
$ws.Cells.Item(17, 8).Value = 18.60327248251037
$ws.Cells.Item(17, 10).Value = 11.20668524117498
$ws.Cells.Item(17, 12).Value = 11.95002347817498
$ws.Cells.Item(17, 13).Value = 18.45817923379168
$ws.Cells.Item(17, 14).Value = 20.01617874933113

$ws.Cells.Item(18, 2).Value = 19.9610684382271
$ws.Cells.Item(18, 4).Value = 8.873265300475722
$ws.Cells.Item(18, 5).Value = 15.19510622094544
$ws.Cells.Item(18, 6).Value = 39.70775668995262
$ws.Cells.Item(18, 7).Value = 49.23535310730666
$ws.Cells.Item(18, 8).Value = 18.59724466983915
$ws.Cells.Item(18, 10).Value = 11.19299403758814
$ws.Cells.Item(18, 12).Value = 11.91178563341997
$ws.Cells.Item(18, 13).Value = 18.42300723794335
$ws.Cells.Item(18, 14).Value = 20.03108100360767

$ws.Cells.Item(19, 2).Value = 19.94150788421073
$ws.Cells.Item(19, 4).Value = 8.866023300949946
$ws.Cells.Item(19, 5).Value = 15.18189613727185
$ws.Cells.Item(19, 6).Value = 39.69861378767495
$ws.Cells.Item(19, 7).Value = 49.19174755417636
$ws.Cells.Item(19, 8).Value = 18.59528361463643
$ws.Cells.Item(19, 10).Value = 11.18835023516113
$ws.Cells.Item(19, 12).Value = 11.89884140242931
$ws.Cells.Item(19, 13).Value = 18.41113230094666
$ws.Cells.Item(19, 14).Value = 20.03616288391847

$ws.Cells.Item(20, 2).Value = 20.02962325550923
$ws.Cells.Item(20, 4).Value = 8.898519261781445
$ws.Cells.Item(20, 5).Value = 15.2411885789802
$ws.Cells.Item(20, 6).Value = 39.74053770200454
$ws.Cells.Item(20, 7).Value = 49.38886008618764
$ws.Cells.Item(20, 8).Value = 18.60442601795707
$ws.Cells.Item(20, 10).Value = 11.20921525562742
$ws.Cells.Item(20, 12).Value = 11.95710134110223
$ws.Cells.Item(20, 13).Value = 18.46470458601121
$ws.Cells.Item(20, 14).Value = 20.01343787827585

$ws.Cells.Item(21, 2).Value = 20.32702715447085
$ws.Cells.Item(21, 4).Value = 9.006029653607921
$ws.Cells.Item(21, 5).Value = 15.43762868739173
$ws.Cells.Item(21, 6).Value = 39.89498871185772
$ws.Cells.Item(21, 7).Value = 50.06600845530978
$ws.Cells.Item(21, 8).Value = 18.64069132725965
$ws.Cells.Item(21, 10).Value = 11.27873256862672
$ws.Cells.Item(21, 12).Value = 12.15285726578073
$ws.Cells.Item(21, 13).Value = 18.64692249620092
$ws.Cells.Item(21, 14).Value = 19.93959999882994

$ws.Cells.Item(22, 2).Value = 20.52211492875835
$ws.Cells.Item(22, 4).Value = 9.075013090440329
$ws.Cells.Item(22, 5).Value = 15.56385976502026
$ws.Cells.Item(22, 6).Value = 40.00598728016975
$ws.Cells.Item(22, 7).Value = 50.51895086934979
$ws.Cells.Item(22, 8).Value = 18.66851694248908
$ws.Cells.Item(22, 10).Value = 11.32371018985877
$ws.Cells.Item(22, 12).Value = 12.28061736220948
$ws.Cells.Item(22, 13).Value = 18.76752449856406
$ws.Cells.Item(22, 14).Value = 19.89321075363269

$ws.Cells.Item(23, 2).Value = 20.41794566208782
$ws.Cells.Item(23, 4).Value = 9.038314740481116
$ws.Cells.Item(23, 5).Value = 15.49669039243219
$ws.Cells.Item(23, 6).Value = 39.94583038463131
$ws.Cells.Item(23, 7).Value = 50.27630215064372
$ws.Cells.Item(23, 8).Value = 18.65328935846546
$ws.Cells.Item(23, 10).Value = 11.2997477263826
$ws.Cells.Item(23, 12).Value = 12.2124574366943
$ws.Cells.Item(23, 13).Value = 18.70302790298409
$ws.Cells.Item(23, 14).Value = 19.91779904447577

$ws.Cells.Item(24, 2).Value = 20.02477857075554
$ws.Cells.Item(24, 4).Value = 8.896740937168769
$ws.Cells.Item(24, 5).Value = 15.23794273613374
$ws.Cells.Item(24, 6).Value = 39.73818413393546
$ws.Cells.Item(24, 7).Value = 49.3779778762005
$ws.Cells.Item(24, 8).Value = 18.60390306604759
$ws.Cells.Item(24, 10).Value = 11.20807160574186
$ws.Cells.Item(24, 12).Value = 11.95390146627444
$ws.Cells.Item(24, 13).Value = 18.46175392616948
$ws.Cells.Item(24, 14).Value = 20.01467634883837

$ws.Cells.Item(25, 2).Value = 19.60563270021157
$ws.Cells.Item(25, 4).Value = 8.738760124038885
$ws.Cells.Item(25, 5).Value = 14.95015732858227
$ws.Cells.Item(25, 6).Value = 39.55782462484831
$ws.Cells.Item(25, 7).Value = 48.458073659822
$ws.Cells.Item(25, 8).Value = 18.56839870580838
$ws.Cells.Item(25, 10).Value = 11.10734913095994
$ws.Cells.Item(25, 12).Value = 11.67545440523733
$ws.Cells.Item(25, 13).Value = 18.20892169817279
$ws.Cells.Item(25, 14).Value = 20.12723670608856
